# "Generate Report for Archive"
# Refresh the localization status report:
#   - the handoff status text moves from "Ready for handoff" to "In Translation"
#     on the Overview sheet (zh-cn/de-de status columns) and on each per-locale
#     sheet's "Status" column.
#   - the Status column(s) are re-sized (narrower, now that the new status text
#     is shorter than the old one) on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status values -------------------------------------------------
# Overview!E2 / F2 hold the per-language ("zh-cn" / "de-de") handoff status.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Each locale sheet keeps the same status in its "Status" column (C2).
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Resize the now-narrower status columns ------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5   # Overview column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # Overview column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5       # zh-cn!C  "Status"
$dede.Columns.Item(3).ColumnWidth = 12.5       # de-de!C  "Status"
